$wb = $excel.ActiveWorkbook

# Fix capitalization of the date/time header labels on both data sheets.
# These cells previously referenced lower-cased shared strings
# ("Start date", "Start day", "Start time", "End time"); they should now
# use the properly capitalized versions ("Start Date", "Start Day",
# "Start Time", "End Time").
$sheetNames = @("ShareSkill", "ManageListings")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("H1").Value = "Start Date"
    $ws.Range("J1").Value = "Start Day"
    $ws.Range("K1").Value = "Start Time"
    $ws.Range("L1").Value = "End Time"
}

# Update the ManageListings sheet selection (no longer the active tab).
$wsManage = $wb.Worksheets.Item("ManageListings")
$wsManage.Range("H1:Q1").Select()

# Make ShareSkill the active sheet/tab with its own selection.
$wsShare = $wb.Worksheets.Item("ShareSkill")
$wsShare.Activate()
$wsShare.Range("E16").Select()
